$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph ending in "...l'address." and insert the new
#    block of paragraphs right after it (before the page-break paragraph).
# ------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*nel campo from specifico di volta in volta l*address.*") {
        $anchorIndex = $i
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchor = $d.Paragraphs($anchorIndex)
$anchor.Range.InsertParagraphAfter()
$placeholder = $d.Paragraphs($anchorIndex + 1)
$placeholderRange = $placeholder.Range

$newXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r><w:t>Flusso log</w:t></w:r>
            <w:r><w:t>ico (setup)</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>
            <w:r><w:t xml:space="preserve"> fatto</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>Mi creo una serie di issuers che si fidano gli uni degli altri (3/4)</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>Il primo fa il createDid() e poi gli altri me li creo con il createChildTrusted.</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>Già pronti all’avvio</w:t></w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r><w:t>Alla verifica:</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r><w:t>Richiedere la Verifiable Credential all’accesso del contenuto</w:t></w:r>
            <w:r><w:t xml:space="preserve"> (rilasciata dall’issuer)</w:t></w:r>
            <w:r><w:t xml:space="preserve">; tramite un server tramite semplice endpoint in locale </w:t></w:r>
            <w:r><w:t>e poi ottenerla</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>
            <w:r><w:t xml:space="preserve"> fatto</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">L’utente in locale si deve generare la VP con ZKP all’interno (dipende dall’età). </w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r><w:t>Il cinema verifica la proof e la non revocation (dipende dal tipo della firma e delle credenziali per generare i dati)</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr>
                <w:ilvl w:val="1"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r><w:t>Nella fase di verifica della proof, ho la verifica dell’issuer (qui il resolveChain)</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$placeholderRange.InsertXML($newXml)

Write-Output "Inserted new block after paragraph $anchorIndex"

# ------------------------------------------------------------------
# 2. The paragraph that holds the page-break run now needs a left
#    indent of 18pt (360 twips) added to its paragraph properties.
#    It is identified as the paragraph immediately preceding "ZKP.".
# ------------------------------------------------------------------
$zkpIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "ZKP.*") {
        $zkpIndex = $i
        break
    }
}
if ($zkpIndex -eq 0) {
    throw "ZKP paragraph not found"
}

$pageBreakPara = $d.Paragraphs($zkpIndex - 1)

# ------------------------------------------------------------------
# 3. Insert a new empty paragraph between the page-break paragraph and
#    the "ZKP." paragraph (done before the indent change below so the
#    new empty paragraph does not inherit the indent).
# ------------------------------------------------------------------
$pageBreakPara.Range.InsertParagraphAfter()
$pageBreakPara.LeftIndent = 18

# ------------------------------------------------------------------
# 4. Strip the lastRenderedPageBreak marker from the "ZKP." paragraph
#    by re-typing its text (this clears out all runs/fields, including
#    the stray lastRenderedPageBreak element, and replaces them with a
#    single plain text run).
# ------------------------------------------------------------------
$zkpIndex2 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "ZKP.*") {
        $zkpIndex2 = $i
        break
    }
}
if ($zkpIndex2 -eq 0) {
    throw "ZKP paragraph not found (pass 2)"
}
$zkpPara = $d.Paragraphs($zkpIndex2)
$zkpPara.Range.Text = "ZKP."

Write-Output "Done"
